$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper approach: for numeric-looking text values in column D, we must force
# Excel to keep them as Text (not auto-convert to a number) while preserving
# the original (default) cell style.

$ws.Range('D2').Value = '49.860.01'
$ws.Range('D3').Value = '2.645.74'
$ws.Range('E3').Value = '  +5.99%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  +0.01%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.75'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +7.62%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '327.21'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +2.31%  '
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('E9').Value = '  +2.86%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.98'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +5.71%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +0.73%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0820'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('E13').Value = '  +1.04%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.36'
$ws.Range('D14').Style = $origStyle
$ws.Range('D15').Value = '3.059.26'
$ws.Range('E15').Value = '  +5.85%  '
$ws.Range('D16').Value = '2.662.57'
$ws.Range('E16').Value = '  +6.51%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.871'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +4.89%  '
$ws.Range('D18').Value = '49.785.57'
$ws.Range('E18').Value = '  +3.90%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.11'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('E20').Value = '  +2.13%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.92'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('D22').Value = '0.0₃0955'
$ws.Range('E22').Value = '  +2.58%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.03'
$ws.Range('D23').Style = $origStyle
$ws.Range('E24').Value = '  +2.20%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.59'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +3.06%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.74'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +3.82%  '
$ws.Range('E27').Value = '  -0.01%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.97'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +2.55%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.71%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.07'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +3.84%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.141'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +0.10%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.20'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('E33').Value = '  +2.74%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.51'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +2.31%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0808'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +4.65%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  +7.10%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.84'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +5.56%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.10'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +8.26%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '125.16'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +2.61%  '
$ws.Range('E41').Value = '  +1.88%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +0.88%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.11'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('E44').Value = '  +3.98%  '
$ws.Range('D45').Value = '2.079.20'
$ws.Range('E45').Value = '  +4.04%  '
$ws.Range('E46').Value = '  +5.54%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.34'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +17.10%  '
$ws.Range('E48').Value = '  +4.95%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +2.47%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.38'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +4.28%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.62'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +4.98%  '
